# update scripts wuth new tpm
# Refresh NATMI ligand/receptor TPM-derived statistics (L1cam-Egfr) with
# newly recomputed values for columns G:J (ligand expr/specificity),
# M:P (receptor expr/specificity) and Q:T (edge weights/specificity)
# across data rows 2-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 3.685507
$ws.Cells.Item(2, 8).Value = 11.056521
$ws.Cells.Item(2, 9).Value = 0.3585631737883472
$ws.Cells.Item(2, 10).Value = 0.3585631737883472
$ws.Cells.Item(2, 13).Value = 0.428743
$ws.Cells.Item(2, 14).Value = 1.286229
$ws.Cells.Item(2, 15).Value = 0.00412050394863168
$ws.Cells.Item(2, 16).Value = 0.00412050394863168
$ws.Cells.Item(2, 17).Value = 1.580135327701
$ws.Cells.Item(2, 18).Value = 14.221217949309
$ws.Cells.Item(2, 19).Value = 0.001477460973428792
$ws.Cells.Item(2, 20).Value = 0.001477460973428792

$ws.Cells.Item(3, 7).Value = 3.685507
$ws.Cells.Item(3, 8).Value = 11.056521
$ws.Cells.Item(3, 9).Value = 0.3585631737883472
$ws.Cells.Item(3, 10).Value = 0.3585631737883472
$ws.Cells.Item(3, 13).Value = 80.22623699999998
$ws.Cells.Item(3, 15).Value = 0.7710272268990069
$ws.Cells.Item(3, 16).Value = 0.7710272268990069
$ws.Cells.Item(3, 17).Value = 295.6743580471589
$ws.Cells.Item(3, 18).Value = 2661.069222424431
$ws.Cells.Item(3, 19).Value = 0.276461969554136
$ws.Cells.Item(3, 20).Value = 0.276461969554136

$ws.Cells.Item(4, 7).Value = 3.685507
$ws.Cells.Item(4, 8).Value = 11.056521
$ws.Cells.Item(4, 9).Value = 0.3585631737883472
$ws.Cells.Item(4, 10).Value = 0.3585631737883472
$ws.Cells.Item(4, 13).Value = 23.39612766666667
$ws.Cells.Item(4, 14).Value = 70.188383
$ws.Cells.Item(4, 15).Value = 0.2248522691523614
$ws.Cells.Item(4, 16).Value = 0.2248522691523614
$ws.Cells.Item(4, 17).Value = 86.22659228839368
$ws.Cells.Item(4, 18).Value = 776.039330595543
$ws.Cells.Item(4, 19).Value = 0.08062374326078239
$ws.Cells.Item(4, 20).Value = 0.08062374326078239

$ws.Cells.Item(5, 9).Value = 0.009647184430711629
$ws.Cells.Item(5, 10).Value = 0.009647184430711629
$ws.Cells.Item(5, 13).Value = 0.428743
$ws.Cells.Item(5, 14).Value = 1.286229
$ws.Cells.Item(5, 15).Value = 0.00412050394863168
$ws.Cells.Item(5, 16).Value = 0.00412050394863168
$ws.Cells.Item(5, 17).Value = 0.042513727137
$ws.Cells.Item(5, 18).Value = 0.382623544233
$ws.Cells.Item(5, 19).Value = 0.00003975126153992534
$ws.Cells.Item(5, 20).Value = 0.00003975126153992534

$ws.Cells.Item(6, 9).Value = 0.009647184430711629
$ws.Cells.Item(6, 10).Value = 0.009647184430711629
$ws.Cells.Item(6, 13).Value = 80.22623699999998
$ws.Cells.Item(6, 15).Value = 0.7710272268990069
$ws.Cells.Item(6, 16).Value = 0.7710272268990069
$ws.Cells.Item(6, 17).Value = 7.955153434682998
$ws.Cells.Item(6, 18).Value = 71.59638091214698
$ws.Cells.Item(6, 19).Value = 0.007438241858994862
$ws.Cells.Item(6, 20).Value = 0.007438241858994862

$ws.Cells.Item(7, 9).Value = 0.009647184430711629
$ws.Cells.Item(7, 10).Value = 0.009647184430711629
$ws.Cells.Item(7, 13).Value = 23.39612766666667
$ws.Cells.Item(7, 14).Value = 70.188383
$ws.Cells.Item(7, 15).Value = 0.2248522691523614
$ws.Cells.Item(7, 16).Value = 0.2248522691523614
$ws.Cells.Item(7, 17).Value = 2.319936623299
$ws.Cells.Item(7, 18).Value = 20.879429609691
$ws.Cells.Item(7, 19).Value = 0.002169191310176842
$ws.Cells.Item(7, 20).Value = 0.002169191310176842

$ws.Cells.Item(8, 7).Value = 6.493877
$ws.Cells.Item(8, 8).Value = 19.481631
$ws.Cells.Item(8, 9).Value = 0.6317896417809412
$ws.Cells.Item(8, 10).Value = 0.6317896417809411
$ws.Cells.Item(8, 13).Value = 0.428743
$ws.Cells.Item(8, 14).Value = 1.286229
$ws.Cells.Item(8, 15).Value = 0.00412050394863168
$ws.Cells.Item(8, 16).Value = 0.00412050394863168
$ws.Cells.Item(8, 17).Value = 2.784204306611
$ws.Cells.Item(8, 18).Value = 25.057838759499
$ws.Cells.Item(8, 19).Value = 0.002603291713662963
$ws.Cells.Item(8, 20).Value = 0.002603291713662963

$ws.Cells.Item(9, 7).Value = 6.493877
$ws.Cells.Item(9, 8).Value = 19.481631
$ws.Cells.Item(9, 9).Value = 0.6317896417809412
$ws.Cells.Item(9, 10).Value = 0.6317896417809411
$ws.Cells.Item(9, 13).Value = 80.22623699999998
$ws.Cells.Item(9, 15).Value = 0.7710272268990069
$ws.Cells.Item(9, 16).Value = 0.7710272268990069
$ws.Cells.Item(9, 17).Value = 520.9793152508489
$ws.Cells.Item(9, 18).Value = 4688.81383725764
$ws.Cells.Item(9, 19).Value = 0.4871270154858761
$ws.Cells.Item(9, 20).Value = 0.487127015485876

$ws.Cells.Item(10, 7).Value = 6.493877
$ws.Cells.Item(10, 8).Value = 19.481631
$ws.Cells.Item(10, 9).Value = 0.6317896417809412
$ws.Cells.Item(10, 10).Value = 0.6317896417809411
$ws.Cells.Item(10, 13).Value = 23.39612766666667
$ws.Cells.Item(10, 14).Value = 70.188383
$ws.Cells.Item(10, 15).Value = 0.2248522691523614
$ws.Cells.Item(10, 16).Value = 0.2248522691523614
$ws.Cells.Item(10, 17).Value = 151.9315753436304
$ws.Cells.Item(10, 18).Value = 1367.384178092673
$ws.Cells.Item(10, 19).Value = 0.1420593345814022
$ws.Cells.Item(10, 20).Value = 0.1420593345814022
